$wb = $excel.ActiveWorkbook

# A new handoff run regenerated the localized artifacts under a fresh
# GUID-based identifier / content hash. Refresh every cell that echoes the
# old identifiers, hashes and timestamps with the newly generated ones.
$oldId = "24f52b42-c635-4002-bb94-6ee3092dcdc8"
$newId = "4cc044c5-93a8-4f3c-9be7-ec7391096bff"
$oldHash = "e9da9db13ca5183eec5eed17937d9da88ca8a992"
$newHash = "8b3f76dc52855dfb4c037a99b9d7d787ea63238a"

$newFileName = "$newId.md"
$newPathAndName = "e2e\$newId.md"
$newGenerateDate = "2016-08-26 15:11:05"
$newZhHandoffFile = "$newId.$newHash.zh-cn.xlf"
$newZhHandoffDate = "2016-08-26 15:10:56"
$newDeHandoffFile = "$newId.$newHash.de-de.xlf"

# The hyperlinks still point at the original (unchanged) GitHub blob URL;
# only the human readable display text advances to the new file name.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/76a8b3da71f8e7681c0fd65f19a7831c2037a9a9/e2e/$oldId.md"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("G2").Value = $newGenerateDate
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkAddress, "", "", $newPathAndName)

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("G2").Value = $newZhHandoffFile
$wsZhCn.Range("H2").Value = $newZhHandoffDate
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkAddress, "", "", $newFileName)

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("G2").Value = $newDeHandoffFile
$wsDeDe.Range("H2").Value = $newGenerateDate
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkAddress, "", "", $newFileName)
